$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before I for "Golongan UH LN" (shifts old I/J/K -> J/K/L)
$ws.Columns("I").Insert()
$ws.Range("I1").Value = "Golongan UH LN"

# Bring the bordered row formatting down into the two new data rows (5 & 6)
$ws.Range("A4:L4").Copy()
$ws.Range("A5:L6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# C5/C6 should not carry the quotePrefix format that C4 has (C4 is quote-prefixed text)
$ws.Range("A4").Copy()
$ws.Range("C5:C6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
# L6 holds a real number (not quote-prefixed text like L4/L5), match a plain bordered cell
$ws.Range("L6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Quote-prefix format (matches the NIK/NPWP/account-number columns) for the new text-numbers
$ws.Range("C2").Copy()
$ws.Range("D5:E6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("L5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 5: Pak Menteri Fulan
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Pak Menteri Fulan"
$ws.Range("D5").Value = "'1234567899877777"
$ws.Range("E5").Value = "'1234567899877777"
$ws.Range("G5").Value = "Menteri"
$ws.Range("I5").Value = "A"
$ws.Range("J5").Value = "Pak Menteri Fulan"
$ws.Range("K5").Value = "BNI"
$ws.Range("L5").Value = "'30302"

# Row 6: Pak Golongan B
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Pak Golongan B"
$ws.Range("D6").Value = "'1234567899877799"
$ws.Range("E6").Value = "'1234567899877799"
$ws.Range("G6").Value = "Duta Besar"
$ws.Range("I6").Value = "B"
$ws.Range("J6").Value = "Pak Golongan B"
$ws.Range("K6").Value = "BNI"
$ws.Range("L6").Value = 31111

# The Golongan/Ruang validation list moved from column O to column P
$dv = $ws.Range("F1").Validation
$dv.Formula1 = "=`$P`$2:`$P`$18"

# Selection left where the author's last edit landed
$ws.Range("I13").Select() | Out-Null
